$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "params": D30 loses its MAX(I30,K30) formula, becomes a literal 2.
# ---------------------------------------------------------------------------
$params = $wb.Worksheets.Item("params")
$params.Range("D30").Value = 2

# ---------------------------------------------------------------------------
# Sheet "levers": rows 2-3 get new formulas/values, rows 4-11 are new data
# following the "C.<capacity>-.<share>" lever-code pattern (Sterman model).
# ---------------------------------------------------------------------------
$levers = $wb.Worksheets.Item("levers")

# Row 2
$levers.Range("A2").Value = 1
$levers.Range("B2").Formula = '="C."&C2&"-."&D2'
$levers.Range("C2").Value = 1
$levers.Range("D2").Value = 0.5

# Row 3
$levers.Range("A3").Formula = "=A2+1"
$levers.Range("B3").Formula = '="C."&C3&"-."&D3'
$levers.Range("C3").Value = 1
$levers.Range("D3").Formula = "=D2+0.1"

# Row 4
$levers.Range("A4").Formula = "=A3+1"
$levers.Range("B4").Formula = '="C."&C4&"-."&D4'
$levers.Range("C4").Value = 1
$levers.Range("D4").Formula = "=D3+0.1"

# Row 5
$levers.Range("A5").Formula = "=A4+1"
$levers.Range("B5").Formula = '="C."&C5&"-."&D5'
$levers.Range("C5").Value = 1
$levers.Range("D5").Formula = "=D4+0.1"

# Row 6
$levers.Range("A6").Formula = "=A5+1"
$levers.Range("B6").Formula = '="C."&C6&"-."&D6'
$levers.Range("C6").Value = 1
$levers.Range("D6").Formula = "=D5+0.1"

# Row 7
$levers.Range("A7").Formula = "=A6+1"
$levers.Range("B7").Formula = '="C."&C7&"-."&D7'
$levers.Range("C7").Value = 2
$levers.Range("D7").Formula = "=D2"

# Row 8
$levers.Range("A8").Formula = "=A7+1"
$levers.Range("B8").Formula = '="C."&C8&"-."&D8'
$levers.Range("C8").Value = 2
$levers.Range("D8").Formula = "=D3"

# Row 9
$levers.Range("A9").Formula = "=A8+1"
$levers.Range("B9").Formula = '="C."&C9&"-."&D9'
$levers.Range("C9").Value = 2
$levers.Range("D9").Formula = "=D4"

# Row 10
$levers.Range("A10").Formula = "=A9+1"
$levers.Range("B10").Formula = '="C."&C10&"-."&D10'
$levers.Range("C10").Value = 2
$levers.Range("D10").Formula = "=D5"

# Row 11
$levers.Range("A11").Value = 10
$levers.Range("B11").Formula = '="C."&C11&"-."&D11'
$levers.Range("C11").Value = 2
$levers.Range("D11").Formula = "=D6"

# ---------------------------------------------------------------------------
# View state: zoom every sheet to 110%, move the active window/selection.
# "levers" becomes the selected/active tab (was "params").
# ---------------------------------------------------------------------------
$configs = $wb.Worksheets.Item("configs")
$varNames = $wb.Worksheets.Item("VariableNames")

$params.Activate()
$params.Range("A30").Select()
$excel.ActiveWindow.Zoom = 110

$configs.Activate()
$configs.Range("D18").Select()
$excel.ActiveWindow.Zoom = 110

$varNames.Activate()
$varNames.Range("D18").Select()
$excel.ActiveWindow.Zoom = 110

$levers.Activate()
$levers.Range("C12").Select()
$excel.ActiveWindow.Zoom = 110
